$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - bump "想去人数" (interest count) values
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1121
$wsExpo.Range("F5").Value = 1009
$wsExpo.Range("F7").Value = 553
$wsExpo.Range("F8").Value = 1192
$wsExpo.Range("F10").Value = 12
$wsExpo.Range("F12").Value = 291
$wsExpo.Range("F14").Value = 90
$wsExpo.Range("F22").Value = 667
$wsExpo.Range("F23").Value = 32
$wsExpo.Range("F24").Value = 641
$wsExpo.Range("F26").Value = 36

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - same counter bumps (rows offset by +1 vs 展览)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1121
$wsAll.Range("F6").Value = 1009
$wsAll.Range("F8").Value = 553
$wsAll.Range("F9").Value = 1192
$wsAll.Range("F12").Value = 12
$wsAll.Range("F14").Value = 291
$wsAll.Range("F16").Value = 90
$wsAll.Range("F30").Value = 667
$wsAll.Range("F31").Value = 32
$wsAll.Range("F32").Value = 641
$wsAll.Range("F34").Value = 36

# ---------------------------------------------------------------------------
# Sheet "全部类型" - a new event row is inserted right after row 43
# (广州·KANAKO ITO&AYANE 2024 LIVE), duplicating that row's data and pushing
# every later row (44-48) down by one (45-49).
# ---------------------------------------------------------------------------
$wsAll.Rows.Item(43).Copy()
$wsAll.Rows.Item(44).Insert()

# The inserted row keeps the running index ("A" column) that the row it
# replaced used to have (43), rather than the copied value (42).
$wsAll.Range("A44").Value = 43
$wsAll.Range("A44").Borders.LineStyle = 1
